$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Ark1")
$ws2 = $wb.Worksheets.Item("Ark2")

# --- Ark1 ("Tidsregistrering") : add missing time-log rows ---------------

# Row 24 gained a Role (column B) entry.
$ws1.Range("B24").Value = "Software Architect"

# Row 25: new task "OC0201,02,03,04" / role "System Analyst " worked
# from 08:50 to 10:30 on 2020-03-09 (serial 43899).
$ws1.Range("A25").Value = "OC0201,02,03,04"
$ws1.Range("B25").Value = "System Analyst "
$ws1.Range("C25").Value = 43899
$ws1.Range("D25").Value = 0.36805555555555558
$ws1.Range("E25").Value = 0.4375

# Row 26: new task "review SD10" worked from 12:20 to 12:40.
$ws1.Range("A26").Value = "review SD10"
$ws1.Range("C26").Value = 43899
$ws1.Range("D26").Value = 0.51388888888888895
$ws1.Range("E26").Value = 0.52777777777777779

# Row 27: new task "Lav SD10" worked from 13:00 to 14:00.
$ws1.Range("A27").Value = "Lav SD10"
$ws1.Range("C27").Value = 43899
$ws1.Range("D27").Value = 0.54166666666666663
$ws1.Range("E27").Value = 0.58333333333333337

# --- View-state tweaks recorded by Excel on save -------------------------

$ws1.Activate()
$ws1.Range("A28").Select()

$ws2.Activate()
$ws2.Range("B5").Select()

$ws1.Activate()
